# Update dSF (column F) values per repull of data / mean calculation fix
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = 9
    3  = -7
    4  = 10
    7  = -9
    9  = -2
    14 = -3
    15 = -3
    17 = 0
    18 = -2
    19 = -5
    23 = 2
    29 = 4
    34 = -1
    37 = 0
    39 = -3
    40 = -4
    42 = -4
    46 = -1
    48 = 2
    52 = -5
    53 = -1
    55 = -2
    61 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

$wb.Save()
